$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 35 and 36 (existing rows 35.. shift down to 37..)
$ws.Rows("35:36").Insert()

# Row 35: new Papaya "Primera" record (Femacal de La Calera, Coquimbo)
$ws.Range("A35").Value = 3
$ws.Range("B35").Value = "Femacal de La Calera"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = 45001
$ws.Range("E35").Value = 5
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100108
$ws.Range("H35").Value = "Tropicales y subtropicales"
$ws.Range("I35").Value = 100108004
$ws.Range("J35").Value = "Papaya"
$ws.Range("K35").Value = "Cultivar IV Región"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 56
$ws.Range("N35").Value = 20000
$ws.Range("O35").Value = 20000
$ws.Range("P35").Value = 20000
$ws.Range("Q35").Value = "$/bandeja 10 kilos"
$ws.Range("R35").Value = "Provincia del Elquí"
$ws.Range("S35").Value = 2000
$ws.Range("T35").Value = 10

# Row 36: new Papaya "Segunda" record (Femacal de La Calera, Coquimbo)
$ws.Range("A36").Value = 3
$ws.Range("B36").Value = "Femacal de La Calera"
$ws.Range("C36").Value = "Coquimbo"
$ws.Range("D36").Value = 45001
$ws.Range("E36").Value = 5
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100108
$ws.Range("H36").Value = "Tropicales y subtropicales"
$ws.Range("I36").Value = 100108004
$ws.Range("J36").Value = "Papaya"
$ws.Range("K36").Value = "Cultivar IV Región"
$ws.Range("L36").Value = "Segunda"
$ws.Range("M36").Value = 54
$ws.Range("N36").Value = 17000
$ws.Range("O36").Value = 17000
$ws.Range("P36").Value = 17000
$ws.Range("Q36").Value = "$/bandeja 10 kilos"
$ws.Range("R36").Value = "Provincia del Elquí"
$ws.Range("S36").Value = 1700
$ws.Range("T36").Value = 10
